$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.352819442749023
$ws.Range("B1").Value = 4.77101993560791
$ws.Range("C1").Value = 2.745707273483276
$ws.Range("D1").Value = 2.393409967422485
$ws.Range("E1").Value = 2.26130485534668
